$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new effort-log entry as row 27:
#   A27 : date 2012-10-17 (serial 41199), using the same date style as the rows above
#   B27 : 1.75 hours of effort
#   D27 : comment, reusing the existing shared string "Manual continued"
$ws.Range("A27").Value = 41199

$ws.Range("B27").Value = 1.75

$ws.Range("D27").Value = "Manual continued"

# Move the active selection to the next empty cell in the row, as Excel
# would leave it after finishing data entry on row 27
$ws.Range("E27").Select()
